# Apply the edits described by the diff between before.xlsx and after.xlsx.
#
# Sheet "展览" (Exhibitions, sheet index 1): "want to go" counters (column F)
# incremented for a number of rows (no rows added/removed).
#
# Sheet "演出" (Performances, sheet index 2): the first data row (old row 2,
# the "杭州·《卡农》永恒经典名曲音乐会" event) is removed entirely; every
# subsequent row shifts up by one. After the shift, the "want to go" count
# for the event that is now in row 2 (previously row 3, the "2024CJMF·不止
# 音乐节" event) increases from 377 to 380.
#
# Sheet "本地生活" (Local life, sheet index 3): "want to go" counters
# incremented for all 4 rows.
#
# Sheet "全部类型" (All types, sheet index 4): "want to go" counters
# incremented for the corresponding rows (no rows added/removed here).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 16
$ws1.Range("F5").Value = 6935
$ws1.Range("F6").Value = 1850
$ws1.Range("F7").Value = 6430
$ws1.Range("F8").Value = 149
$ws1.Range("F9").Value = 1991
$ws1.Range("F10").Value = 534
$ws1.Range("F11").Value = 28
$ws1.Range("F12").Value = 14
$ws1.Range("F17").Value = 8223
$ws1.Range("F18").Value = 148
$ws1.Range("F22").Value = 1779
$ws1.Range("F23").Value = 854
$ws1.Range("F28").Value = 179
$ws1.Range("F29").Value = 5
$ws1.Range("F30").Value = 1918
$ws1.Range("F31").Value = 828
$ws1.Range("F32").Value = 425
$ws1.Range("F33").Value = 4
$ws1.Range("F34").Value = 9
$ws1.Range("F35").Value = 142
$ws1.Range("F36").Value = 117
$ws1.Range("F38").Value = 3942

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# Remove the "《卡农》永恒经典名曲音乐会" row (old row 2); everything below
# shifts up by one row, and the dimension shrinks from A1:I27 to A1:I26.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()

# The row that is now row 2 (previously row 3) gets an updated "want to go"
# count, from 377 to 380.
$ws2.Range("F2").Value = 380

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9551
$ws3.Range("F3").Value = 2304
$ws3.Range("F4").Value = 698
$ws3.Range("F5").Value = 289

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9551
$ws4.Range("F3").Value = 2304
$ws4.Range("F4").Value = 698
$ws4.Range("F7").Value = 380
$ws4.Range("F8").Value = 6935
$ws4.Range("F10").Value = 1851
$ws4.Range("F11").Value = 6430
$ws4.Range("F12").Value = 149
$ws4.Range("F13").Value = 1991
$ws4.Range("F15").Value = 534
$ws4.Range("F16").Value = 14
$ws4.Range("F23").Value = 8223
$ws4.Range("F24").Value = 148
$ws4.Range("F28").Value = 1779
$ws4.Range("F29").Value = 854
$ws4.Range("F32").Value = 179
$ws4.Range("F33").Value = 1918
$ws4.Range("F34").Value = 828
$ws4.Range("F36").Value = 425
$ws4.Range("F37").Value = 4
$ws4.Range("F39").Value = 9
$ws4.Range("F41").Value = 117
$ws4.Range("F44").Value = 3942
